# Insert a new record row at row 255 (pushing the existing rows 255..364 down
# to 256..365) for "Fruta / hortaliza, semanal" - Femacal de La Calera /
# Arandano (blue).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 255..364 down by one row.
$ws.Rows.Item(255).Insert()

# Populate the newly-inserted row 255 with the new record.
$ws.Cells.Item(255, 1).Value = 3
$ws.Cells.Item(255, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(255, 3).Value = "Coquimbo"
$ws.Cells.Item(255, 4).Value = 45161
$ws.Cells.Item(255, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(255, 5).Value = 5
$ws.Cells.Item(255, 6).Value = "Fruta"
$ws.Cells.Item(255, 7).Value = 100101
$ws.Cells.Item(255, 8).Value = "Berries"
$ws.Cells.Item(255, 9).Value = 100101001
$ws.Cells.Item(255, 10).Value = "Arándano (blue)"
$ws.Cells.Item(255, 11).Value = "Sin especificar"
$ws.Cells.Item(255, 12).Value = "Primera"
$ws.Cells.Item(255, 13).Value = 45
$ws.Cells.Item(255, 14).Value = 14000
$ws.Cells.Item(255, 15).Value = 14000
$ws.Cells.Item(255, 16).Value = 14000
$ws.Cells.Item(255, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(255, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(255, 19).Value = 9333
$ws.Cells.Item(255, 20).Value = 1.5
